# Updated cryptos list on Sun Apr  2 07:31:00 UTC 2023 with GitHub Actions
#
# Refreshes the per-coin "Price" (column D) and "Volume(1h)" (column E)
# snapshot values scraped from coinranking.com. Columns D/E hold plain text
# (not numbers) in this sheet, so for the price cells whose new reading
# happens to look numeric we pin NumberFormat to Text ("@") first -- this
# mirrors how the sheet is normally populated and stops Excel's COM layer
# from silently reinterpreting e.g. "317.19" as a float.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")


$ws.Range("D2").Value = '28.545.49'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '1.828.03'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.19'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5171'
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3890'
$ws.Range("E8").Value = '  -1.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08393'
$ws.Range("E9").Value = '  +8.90%  '
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.96'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.424'
$ws.Range("E12").Value = '  +2.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.32'
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.528'
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").Value = '1.827.06'
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.48'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001130'
$ws.Range("E18").Value = '  +4.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06638'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.78'
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.082'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").Value = '28.578.81'
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.44'
$ws.Range("E24").Value = '  +2.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.278'
$ws.Range("E25").Value = '  +1.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.18'
$ws.Range("E26").Value = '  +2.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.67'
$ws.Range("E27").Value = '  +1.63%  '
$ws.Range("D28").Value = '2.034.98'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.419'
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.76'
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("E32").Value = '  -2.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.744'
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07586'
$ws.Range("E34").Value = '  +5.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.674'
$ws.Range("E35").Value = '  +1.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2230'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02383'
$ws.Range("E37").Value = '  +2.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.264'
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.767'
$ws.Range("E39").Value = '  -2.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6385'
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.47'
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.194'
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.400'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.59'
$ws.Range("E44").Value = '  +1.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6031'
$ws.Range("E45").Value = '  +2.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.779'
$ws.Range("E46").Value = '  +1.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '127.14'
$ws.Range("E47").Value = '  +2.30%  '
$ws.Range("E48").Value = '  +0.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.207'
$ws.Range("E49").Value = '  +2.16%  '
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.79'
$ws.Range("E51").Value = '  +1.28%  '
